# OIE -> WOAH rename across all worksheets (text content only; URLs containing
# "oie.int" are left untouched since the source diff does not alter them).
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $val = $cell.Value2

            if ($val -ne $null -and $val -is [string] -and $val.Contains("OIE")) {
                # Replace the standalone token "OIE" with "WOAH" everywhere it
                # occurs in the text, but leave URLs (e.g. oie.int, already
                # lower-case) untouched since those only match lower-case "oie".
                $newVal = $val.Replace("OIE", "WOAH")
                $cell.Value2 = $newVal
            }
        }
    }
}
